{"js": "// Merge the split \"Example 10: ... Finding \" runs in the chapter heading\n// paragraph into a single run (keeping the formatting of the first run),\n// and change the \"zscore = qnorm(0.975)\" assignment operator from \"=\" to\n// the R-style \"<-\" assignment, preserving the existing run/style\n// structure (OtherTok) of that source-code paragraph.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// --- 1) Merge \"Example 10: ... Finding \" runs into one run -------------\nconst headingText = \"Example 10:  Estimating Mean Education in South Africa \u2013 Finding \";\nfor (const para of paragraphs.items) {\n  if (para.text.indexOf(headingText) !== -1) {\n    const found = para.search(headingText, { matchCase: true });\n    found.load(\"items\");\n    await context.sync();\n    if (found.items.length > 0) {\n      found.items[0].insertText(headingText, \"Replace\");\n      await context.sync();\n    }\n    break;\n  }\n}\n\n// --- 2) \"zscore = qnorm(0.975)\" -> \"zscore <- qnorm(0.975)\" ------------\nfor (const para of paragraphs.items) {\n  if (para.text.indexOf(\"zscore =\") !== -1) {\n    const eq = para.search(\"=\", { matchCase: true });\n    eq.load(\"items\");\n    await context.sync();\n    if (eq.items.length > 0) {\n      eq.items[0].insertText(\"<-\", \"Replace\");\n      await context.sync();\n    }\n    break;\n  }\n}\n", "ps1": "# Merge the split \"Example 10: ... Finding \" runs in the chapter heading\n# paragraph into a single run (keeping the formatting of the first run),\n# and change the \"zscore = qnorm(0.975)\" assignment operator from \"=\" to\n# the R-style \"<-\" assignment, preserving the existing run/style\n# structure (OtherTok) of that source-code paragraph.\n\n$d = $word.ActiveDocument\n\n$headingText = \"Example 10:  Estimating Mean Education in South Africa \" + [char]0x2013 + \" Finding \"\n\n# --- 1) Merge \"Example 10: ... Finding \" runs into one run -------------\n# A straight Find+Replace with *identical* text is treated as a no-op by\n# the engine (no run merge happens), so first swap in a unique marker and\n# then replace the marker with the real text - that guarantees an actual\n# content diff that collapses the surrounding runs into one.\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Estimating Mean Education*\") {\n        $r = $p.Range\n        $find = $r.Find\n        $find.Text = $headingText\n        $find.MatchCase = $true\n        $found = $find.Execute()\n        if ($found) {\n            $r.Text = \"@@HEADING_MARKER@@\"\n        }\n        break\n    }\n}\n\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*HEADING_MARKER*\") {\n        $r2 = $p.Range\n        $find2 = $r2.Find\n        $find2.Text = \"@@HEADING_MARKER@@\"\n        $find2.MatchCase = $true\n        $found2 = $find2.Execute()\n        if ($found2) {\n            $r2.Text = $headingText\n        }\n        break\n    }\n}\n\n# --- 2) \"zscore = qnorm(0.975)\" -> \"zscore <- qnorm(0.975)\" ------------\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*zscore =*\") {\n        $r3 = $p.Range\n        $find3 = $r3.Find\n        $find3.Text = \"=\"\n        $find3.MatchCase = $true\n        $found3 = $find3.Execute()\n        if ($found3) {\n            $r3.Text = \"<-\"\n        }\n        break\n    }\n}\n"}
